$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/8/2024  Through  7/14/2024"

# --- Text-valued cells that change from numeric style to shared-text style (full copy, value matches source) ---
$ws.Range("C28").Copy($ws.Range("G28"))
$ws.Range("M28").Copy($ws.Range("H28"))
$ws.Range("C33").Copy($ws.Range("D33"))
$ws.Range("M33").Copy($ws.Range("E33"))

# --- Numeric cells that change style: copy number format from a same-target-style neighbor, then set value ---
$ws.Range("C16").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 1
$ws.Range("H16").Copy($ws.Range("E16"))
$ws.Range("E16").Value = 0
$ws.Range("G18").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 2
$ws.Range("H18").Copy($ws.Range("E18"))
$ws.Range("E18").Value = -100
$ws.Range("F25").Copy($ws.Range("D25"))
$ws.Range("D25").Value = 4
$ws.Range("H25").Copy($ws.Range("E25"))
$ws.Range("E25").Value = -75
$ws.Range("C16").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 1
$ws.Range("H16").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100
$ws.Range("C16").Copy($ws.Range("G29"))
$ws.Range("G29").Value = 1
$ws.Range("H16").Copy($ws.Range("H29"))
$ws.Range("H29").Value = -100
$ws.Range("C16").Copy($ws.Range("J29"))
$ws.Range("J29").Value = 1
$ws.Range("H16").Copy($ws.Range("K29"))
$ws.Range("K29").Value = -100
$ws.Range("C16").Copy($ws.Range("D30"))
$ws.Range("D30").Value = 1
$ws.Range("H16").Copy($ws.Range("E30"))
$ws.Range("E30").Value = -100
$ws.Range("C16").Copy($ws.Range("G30"))
$ws.Range("G30").Value = 1
$ws.Range("H16").Copy($ws.Range("H30"))
$ws.Range("H30").Value = -100
$ws.Range("C16").Copy($ws.Range("J30"))
$ws.Range("J30").Value = 1
$ws.Range("H16").Copy($ws.Range("K30"))
$ws.Range("K30").Value = -100

# --- Same-style numeric value updates ---
$ws.Range("F16").Value = 2
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 15
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = 36.363636363636
$ws.Range("L16").Value = 114.285714285714
$ws.Range("M16").Value = 15.384615384615
$ws.Range("N16").Value = -44.444444444444
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -46.153846153846
$ws.Range("I17").Value = 51
$ws.Range("J17").Value = 54
$ws.Range("K17").Value = -5.555555555555
$ws.Range("L17").Value = 112.5
$ws.Range("M17").Value = 131.818181818182
$ws.Range("N17").Value = -15
$ws.Range("G18").Value = 7
$ws.Range("J18").Value = 32
$ws.Range("K18").Value = -50
$ws.Range("L18").Value = -5.882352941176
$ws.Range("N18").Value = -89.808917197452
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 66.666666666666
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 26
$ws.Range("H19").Value = 19.230769230769
$ws.Range("I19").Value = 168
$ws.Range("J19").Value = 157
$ws.Range("K19").Value = 7.006369426751
$ws.Range("L19").Value = 16.666666666666
$ws.Range("M19").Value = 136.619718309859
$ws.Range("N19").Value = 54.128440366972
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = -50
$ws.Range("J20").Value = 45
$ws.Range("K20").Value = -51.111111111111
$ws.Range("N20").Value = -94.021739130434
$ws.Range("C21").Value = 7
$ws.Range("E21").Value = -46.153846153846
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = -21.818181818181
$ws.Range("I21").Value = 276
$ws.Range("J21").Value = 301
$ws.Range("K21").Value = -8.305647840531
$ws.Range("L21").Value = 14.049586776859
$ws.Range("M21").Value = 45.263157894736
$ws.Range("N21").Value = -61.931034482758
$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -73.684210526315
$ws.Range("F24").Value = 27
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = -34.146341463414
$ws.Range("I24").Value = 211
$ws.Range("J24").Value = 264
$ws.Range("K24").Value = -20.075757575757
$ws.Range("L24").Value = -1.401869158878
$ws.Range("M24").Value = -22.992700729927
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 28.571428571428
$ws.Range("I25").Value = 97
$ws.Range("J25").Value = 115
$ws.Range("K25").Value = -15.652173913043
$ws.Range("L25").Value = 94
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = -18.75
$ws.Range("I26").Value = 89
$ws.Range("J26").Value = 105
$ws.Range("K26").Value = -15.238095238095
$ws.Range("L26").Value = -2.197802197802
$ws.Range("M26").Value = -20.535714285714
$ws.Range("F28").Value = 2
$ws.Range("L33").Value = -80
